# Atualização de bases das ligas, do dia: 15-04-2024 às 22:35
#
# The source feed re-sorted/re-paired several match records between two
# consecutive rows (everything except the running index in column A), and
# one existing match record (row 160) was corrected/expanded with new
# odds data (including previously-missing FT result columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param($ws, [int]$row1, [int]$row2, [int]$colStart, [int]$colEnd)

    # Snapshot both rows first so the swap is not order-dependent.
    $vals1 = @{}
    $vals2 = @{}
    for ($c = $colStart; $c -le $colEnd; $c++) {
        $vals1[$c] = $ws.Cells.Item($row1, $c).Value()
        $vals2[$c] = $ws.Cells.Item($row2, $c).Value()
    }
    for ($c = $colStart; $c -le $colEnd; $c++) {
        $ws.Cells.Item($row1, $c).Value = $vals2[$c]
        $ws.Cells.Item($row2, $c).Value = $vals1[$c]
    }
}

# Row index (column A) stays put; columns B(2) .. AC(29) swap wholesale
# between each of these adjacent-row pairs.
Swap-RowData $ws 58 59 2 29
Swap-RowData $ws 89 90 2 29
Swap-RowData $ws 91 92 2 29
Swap-RowData $ws 103 104 2 29
Swap-RowData $ws 108 109 2 29

# Row 160: id 8088065 -> 8088066, with corrected date/teams/odds and newly
# populated FTHG/FTAG/FTR (H160/I160/J160) plus PL_AhOver/PL_AhUnder
# (AB160/AC160).
$ws.Cells.Item(160, 2).Value = 8088066
$ws.Cells.Item(160, 5).Value = 45396.47916666666
$ws.Cells.Item(160, 6).Value = "Al Shorta SC"
$ws.Cells.Item(160, 7).Value = "Al Minaa"
$ws.Cells.Item(160, 8).Value = 2
$ws.Cells.Item(160, 9).Value = 0
$ws.Cells.Item(160, 10).Value = "H"
$ws.Cells.Item(160, 11).Value = 1.333
$ws.Cells.Item(160, 12).Value = 5
$ws.Cells.Item(160, 13).Value = 6
$ws.Cells.Item(160, 14).Value = 1.25
$ws.Cells.Item(160, 15).Value = 6
$ws.Cells.Item(160, 16).Value = 7
$ws.Cells.Item(160, 17).Value = -1.75
$ws.Cells.Item(160, 21).Value = 1.9
$ws.Cells.Item(160, 22).Value = 1.9
$ws.Cells.Item(160, 23).Value = 0.25
$ws.Cells.Item(160, 24).Value = -1
$ws.Cells.Item(160, 25).Value = -1
$ws.Cells.Item(160, 26).Value = 0.475
$ws.Cells.Item(160, 27).Value = -0.5
$ws.Cells.Item(160, 28).Value = -1
$ws.Cells.Item(160, 29).Value = 0.8999999999999999
